$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text is preserved verbatim (avoid Excel auto-converting numeric-looking
# strings like "0.0000404" or "2.80" into floating point numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.418.69'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +5.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.714.96'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +6.99%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '422.99'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.77'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +2.00%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.708.33'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +7.12%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.06%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.768'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.185'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +11.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000404'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +47.99%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '43.28'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.28'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.305.53'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +7.10%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.83'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.737.04'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +8.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.08'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +5.37%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.97%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '67.448.10'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +6.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '453.09'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.63'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +14.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '89.81'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -1.58%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '38.03'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +11.41%  '
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.34'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Filecoin'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.18'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.52%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.124'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +10.38%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.80'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.77%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.42'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.164'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.23%  '
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '42.05'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +5.21%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '56.54'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.99%  '
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0496'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +1.53%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0785'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +19.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.13'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +34.88%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +3.75%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.41'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.52'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +25.50%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '148.23'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.94%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.44'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.93'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.59%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +5.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.68'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.309'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.64%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +16.25%  '
